$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string used in column B (e.g. "EEE_150_FFF" -> "EEE_150_FFF1")
$ws.Range("B2:B5").Value = "EEE_150_FFF1"

# Widen column B to fit new text (matches Excel's bestFit recalculation after
# the longer string was entered; closest attainable width on this host)
$ws.Columns.Item(2).ColumnWidth = 12.3

# Move the active selection to D12
$ws.Range("D12").Select()
